$d = $word.ActiveDocument

# --- Change 1: merge the "Head of the Schoo" + "l" runs into a single run
#     reading "Head of the School". ---
$d.Content.Find.Execute("Head of the Schoo" + "l", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Head of the School", 2) | Out-Null

# --- Locate the "Website Visitors" (Heading 3) paragraph, and the empty
#     paragraph that immediately precedes it, by content rather than a
#     hard-coded paragraph index. ---
$count = $d.Paragraphs.Count
$visitorsIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("Website Visitors")) {
        $visitorsIndex = $i
        break
    }
}
$noteIndex = $visitorsIndex - 1

# --- Change 2: fill that empty paragraph with the new "NOTE: ..." runs
#     (bold "NOTE: " lead-in carrying the page-break hint, Calibri
#     throughout, split into the same run boundaries as the authored text).
#     The paragraph's own <w:p>/<w:pPr> is read back via WordOpenXML and
#     reused verbatim so paraId/rsid metadata are left exactly as-is. ---
$notePara = $d.Paragraphs.Item($noteIndex)
$noteFullXml = $notePara.Range.WordOpenXML
$noteFullXml -match '<w:body>(<w:p\b[^>]*>(?:<w:pPr>.*?</w:pPr>)?)' | Out-Null
$noteParaPrefix = $matches[1]

$noteRuns = (
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/></w:rPr>' + `
        '<w:lastRenderedPageBreak/><w:t xml:space="preserve">NOTE: </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t xml:space="preserve">It </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t xml:space="preserve">is worth noting </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t xml:space="preserve">that the ''Course Coordinator'' field was omitted from the Courses model. This decision was </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t>since</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t xml:space="preserve"> the official Curtin University course </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t>catalogue</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t xml:space="preserve"> does not provide specific details on course coordinators.</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>' + `
        '<w:t xml:space="preserve"> </w:t></w:r>'
)

$noteXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body>' + $noteParaPrefix + $noteRuns + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$notePara.Range.InsertXML($noteXml)

# --- Change 3: re-fetch "Website Visitors" (its paragraph index hasn't
#     moved - we only added content to the paragraph before it) and
#     rewrite its run without the now-relocated lastRenderedPageBreak,
#     reusing the paragraph's own <w:p>/<w:pPr>/<w:r> opening tags so
#     nothing else about it changes. ---
$visitorsPara = $d.Paragraphs.Item($visitorsIndex)
$visitorsFullXml = $visitorsPara.Range.WordOpenXML
$visitorsFullXml -match '<w:body>(<w:p\b[^>]*>(?:<w:pPr>.*?</w:pPr>)?)(<w:r\b[^>]*>)' | Out-Null
$visitorsParaPrefix = $matches[1]
$visitorsRunOpenTag = $matches[2]

$visitorsXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body>' + $visitorsParaPrefix + $visitorsRunOpenTag + '<w:t>Website Visitors</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$visitorsPara.Range.InsertXML($visitorsXml)
